$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9091022
$ws.Range("I9").Value = 100000000
$ws.Range("J9").Value = 124
$ws.Range("K9").Value = 100000000
$ws.Range("L9").Value = 124
$ws.Range("M9").Value = -99999831
$ws.Range("N9").Value = -462

$ws.Range("H19").Value = 8080400
$ws.Range("I19").Value = 6708591.5
$ws.Range("J19").Value = 10000932
$ws.Range("K19").Value = 6708591.5
$ws.Range("L19").Value = 10000932
$ws.Range("M19").Value = -6708416.5
$ws.Range("N19").Value = -10001282

$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4652

$ws.Range("H43").Value = 999.1875
$ws.Range("I43").Value = 795
$ws.Range("J43").Value = 1092
$ws.Range("K43").Value = 795
$ws.Range("L43").Value = 1092
$ws.Range("M43").Value = -726
$ws.Range("N43").Value = -1230

$ws.Range("H55").Value = 228.92857
$ws.Range("I55").Value = 130.3
$ws.Range("K55").Value = 130.3
$ws.Range("M55").Value = 83.69999999999999

$ws.Range("H112").Value = 1344.3784
$ws.Range("I112").Value = 750
$ws.Range("J112").Value = 1378.3429
$ws.Range("K112").Value = 2250
$ws.Range("L112").Value = 4135.028700000001
$ws.Range("M112").Value = -1142
$ws.Range("N112").Value = -6351.028700000001

$ws.Range("H115").Value = 14286435
$ws.Range("I115").Value = 14286435
$ws.Range("K115").Value = 42859305
$ws.Range("M115").Value = -42857738

$ws.Range("H116").Value = 4052598.8
$ws.Range("I116").Value = 19233270
$ws.Range("J116").Value = 4420
$ws.Range("K116").Value = 19233270
$ws.Range("L116").Value = 4420
$ws.Range("M116").Value = -19229828
$ws.Range("N116").Value = -11304

$ws.Range("H132").Value = 2747.9666
$ws.Range("I132").Value = 2314.7646
$ws.Range("J132").Value = 5202.778
$ws.Range("K132").Value = 6944.293799999999
$ws.Range("L132").Value = 15608.334
$ws.Range("M132").Value = -4414.293799999999
$ws.Range("N132").Value = -20668.334

$ws.Range("H137").Value = 23907.89
$ws.Range("I137").Value = 1369.0344
$ws.Range("J137").Value = 62356.53
$ws.Range("K137").Value = 4107.1032
$ws.Range("L137").Value = 187069.59
$ws.Range("M137").Value = -1557.1032
$ws.Range("N137").Value = -192169.59

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H32").Value = 18523364
$ws.Range("I32").Value = 20836534
$ws.Range("J32").Value = 18002.166
$ws.Range("K32").Value = 20836534
$ws.Range("L32").Value = 18002.166
$ws.Range("M32").Value = -20836247
$ws.Range("N32").Value = -18576.166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1912.5714
$ws.Range("I99").Value = 1297.8572
$ws.Range("J99").Value = 3142
$ws.Range("K99").Value = 1297.8572
$ws.Range("L99").Value = 3142
$ws.Range("M99").Value = 200.1428000000001
$ws.Range("N99").Value = -6138

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 475
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H31").Value = 1936.5111
$ws.Range("I31").Value = 1451.0571
$ws.Range("J31").Value = 3635.6
$ws.Range("K31").Value = 1451.0571
$ws.Range("L31").Value = 3635.6
$ws.Range("M31").Value = -1156.0571
$ws.Range("N31").Value = -4225.6

$ws.Range("H34").Value = 1936.5111
$ws.Range("I34").Value = 1451.0571
$ws.Range("J34").Value = 3635.6
$ws.Range("K34").Value = 1451.0571
$ws.Range("L34").Value = 3635.6
$ws.Range("M34").Value = -1249.0571
$ws.Range("N34").Value = -4039.6

$ws.Range("H58").Value = 784.3484999999999
$ws.Range("I58").Value = 701.2406999999999
$ws.Range("J58").Value = 1158.3334
$ws.Range("K58").Value = 701.2406999999999
$ws.Range("L58").Value = 1158.3334
$ws.Range("M58").Value = -498.2406999999999
$ws.Range("N58").Value = -1564.3334

$ws.Range("H86").Value = 4735.6
$ws.Range("I86").Value = 7017.8335
$ws.Range("J86").Value = 3214.111
$ws.Range("K86").Value = 7017.8335
$ws.Range("L86").Value = 3214.111
$ws.Range("M86").Value = -5894.8335
$ws.Range("N86").Value = -5460.111

$ws.Range("H89").Value = 4735.6
$ws.Range("I89").Value = 7017.8335
$ws.Range("J89").Value = 3214.111
$ws.Range("K89").Value = 35089.1675
$ws.Range("L89").Value = 16070.555
$ws.Range("M89").Value = -29473.1675
$ws.Range("N89").Value = -27302.555

$ws.Range("H107").Value = 553.9524
$ws.Range("I107").Value = 351.92856
$ws.Range("J107").Value = 958
$ws.Range("K107").Value = 351.92856
$ws.Range("L107").Value = 958
$ws.Range("M107").Value = 1568.07144
$ws.Range("N107").Value = -4798

$ws.Range("H136").Value = 784.3484999999999
$ws.Range("I136").Value = 701.2406999999999
$ws.Range("J136").Value = 1158.3334
$ws.Range("K136").Value = 2103.7221
$ws.Range("L136").Value = 3475.0002
$ws.Range("M136").Value = 446.2779
$ws.Range("N136").Value = -8575.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 6500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 19500
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -19948

$ws.Range("H12").Value = 5.8
$ws.Range("I12").Value = 13.5
$ws.Range("K12").Value = 40.5
$ws.Range("M12").Value = 132.5

$ws.Range("H23").Value = 98.57143000000001
$ws.Range("I23").Value = 90
$ws.Range("J23").Value = 105
$ws.Range("K23").Value = 270
$ws.Range("L23").Value = 315
$ws.Range("M23").Value = -35
$ws.Range("N23").Value = -785

$ws.Range("H131").Value = 716.1129
$ws.Range("I131").Value = 414.17392
$ws.Range("J131").Value = 894.1795
$ws.Range("K131").Value = 1242.52176
$ws.Range("L131").Value = 2682.5385
$ws.Range("M131").Value = 3797.47824
$ws.Range("N131").Value = -12762.5385

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2076.75
$ws.Range("I122").Value = 1802.3334
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 5407.0002
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -2957.0002
$ws.Range("N122").Value = -13600

$ws.Range("H126").Value = 3083.0908
$ws.Range("I126").Value = 1979.7142
$ws.Range("J126").Value = 5014
$ws.Range("K126").Value = 5939.142599999999
$ws.Range("L126").Value = 15042
$ws.Range("M126").Value = -3469.142599999999
$ws.Range("N126").Value = -19982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1061.4286
$ws.Range("J46").Value = 1107.5
$ws.Range("L46").Value = 1107.5
$ws.Range("N46").Value = -1483.5

$ws.Range("H55").Value = 687.3333
$ws.Range("I55").Value = 655.8333
$ws.Range("J55").Value = 813.3333
$ws.Range("K55").Value = 655.8333
$ws.Range("L55").Value = 813.3333
$ws.Range("M55").Value = -482.8333
$ws.Range("N55").Value = -1159.3333

$ws.Range("H134").Value = 43329
$ws.Range("J134").Value = 43329
$ws.Range("L134").Value = 43329
$ws.Range("N134").Value = -53469

$ws.Range("H137").Value = 35914.8
$ws.Range("J137").Value = 35914.8
$ws.Range("L137").Value = 35914.8
$ws.Range("N137").Value = -46114.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 10000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -9712
$ws.Range("N15").ClearContents()

$ws.Range("H112").Value = 28257.666
$ws.Range("J112").Value = 28257.666
$ws.Range("L112").Value = 28257.666
$ws.Range("N112").Value = -31211.666
